$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "...por un periodo de {{PERIODO_MESES}} meses, a fin de..."
#    becomes
#    "...por un periodo de {{PERIODO_MESES}}, a fin de..."
#    i.e. the word " meses" (leading space included) right before the
#    comma is removed.
# ------------------------------------------------------------------
$targetParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("complementar la formaci")) {
        $targetParaIndex = $i
        break
    }
}

if ($targetParaIndex -ne $null) {
    $para = $d.Paragraphs.Item($targetParaIndex)
    $paraText = $para.Range.Text
    $cut = " meses,"
    $cutPos = $paraText.IndexOf($cut)
    if ($cutPos -ge 0) {
        $absStart = $para.Range.Start + $cutPos
        # remove just the " meses" part, keep the comma that follows it
        $rng = $d.Range($absStart, $absStart + 6)
        $rng.Text = ""
    }
}

# ------------------------------------------------------------------
# 2. The hidden "_GoBack" bookmark (marks the last edit location) sat at
#    the end of that paragraph; after this edit Word leaves it on the
#    following, empty paragraph instead.
# ------------------------------------------------------------------
if ($targetParaIndex -ne $null) {
    $nextPara = $d.Paragraphs.Item($targetParaIndex + 1)
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
    $d.Bookmarks.Add("_GoBack", $nextPara.Range)
}
